# Commit: "updated scripts to support Nodata (-1) from CSV tables and upgraded to Python 3"
#
# This fills in explicit Nodata sentinel values (-1) in columns C:F for every
# body row that previously had no radii recorded, matching the new CSV-driven
# export. A handful of rows that already had partial data only need a single
# missing column backfilled, and one row (158) had its existing -1 values
# re-typed from a 2-decimal style to an integer style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Contiguous blocks of rows whose C:F cells are completely empty and need
# the full Nodata (-1) block, formatted as integers (no decimals).
# (Kept as parallel start/end arrays rather than an array-of-arrays, since
# nested arrays get flattened by this host's PowerShell engine.)
$blockStarts = @(26, 79,  96, 110, 130, 151, 164, 166, 174)
$blockEnds   = @(59, 91, 108, 113, 141, 155, 164, 166, 176)

for ($i = 0; $i -lt $blockStarts.Length; $i++) {
  $r1 = $blockStarts[$i]
  $r2 = $blockEnds[$i]
  $rng = $ws.Range("C${r1}:F${r2}")
  $rng.Value = -1
  $rng.NumberFormat = "0"
}

# Row 158: E:F already held -1 but used the 2-decimal style; restyle them to
# the integer Nodata style (values are unchanged).
$ws.Range("E158:F158").NumberFormat = "0"

# Row 159: only column C was missing; backfill it with -1, keeping the
# existing 2-decimal style used by the rest of that row.
$ws.Range("C159").Value = -1
$ws.Range("C159").NumberFormat = "0.00"

# Rows 169, 177, 178: only column C was missing; backfill with -1 using the
# integer Nodata style.
$ws.Range("C169").Value = -1
$ws.Range("C169").NumberFormat = "0"

$ws.Range("C177").Value = -1
$ws.Range("C177").NumberFormat = "0"

$ws.Range("C178").Value = -1
$ws.Range("C178").NumberFormat = "0"

# Reset the view: scroll back to the top and move the selection to B2
# (previously the view was scrolled down to row 140 with F159 selected).
[void]$ws.Range("B2").Select()
